$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.941.17"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "2.365.14"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.35"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.18"
$ws.Range("E6").Value = "  -3.59%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  -2.66%  "

$ws.Range("D9").Value = "2.364.55"
$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("E11").Value = "  +1.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.91"
$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").Value = "2.775.63"
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("D17").Value = "60.868.41"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "2.374.60"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.49"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.11"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.93"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.65"
$ws.Range("E25").Value = "  -15.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.22"
$ws.Range("E26").Value = "  +0.81%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "2.474.59"
$ws.Range("E28").Value = "  -1.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  -0.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.148"
$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("D31").Value = "0.0₃0875"
$ws.Range("E31").Value = "  -7.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "495.49"
$ws.Range("E32").Value = "  -8.39%  "

$ws.Range("E33").Value = "  -5.43%  "

$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -5.73%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.377"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.50"
$ws.Range("E40").Value = "  +1.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.30"
$ws.Range("E41").Value = "  -5.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.24"
$ws.Range("E42").Value = "  +4.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.98"
$ws.Range("E44").Value = "  +1.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.53"
$ws.Range("E45").Value = "  +3.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.58"
$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").Value = "  -9.25%  "

$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.572"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("E50").Value = "  -5.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0905"
$ws.Range("E51").Value = "  -0.60%  "

# Reset number format back to General/Normal style so cells match default (no explicit style)
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
